# DataSource - Emision Motor - Inspeccion.xlsx
# Refresh the test data block (rows 7-8): bump the inspection date and the
# Motor/Chasis reference codes to the next batch, and correct the vehicle
# year on row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FechaInicio / inspection date for rows 7 & 8 (shared by both rows)
$ws.Range("K7").Value = "'01/02/2021"
$ws.Range("K8").Value = "'01/02/2021"

# Motor / Chasis reference codes - row 7: RGR005 -> RGR013
$ws.Range("T7").Value = "RGR013"
$ws.Range("U7").Value = "ABCD0RGR013"
$ws.Range("V7").Value = "ZXC0987RGR013"

# Motor / Chasis reference codes - row 8: RGR006 -> RGR014
$ws.Range("T8").Value = "RGR014"
$ws.Range("U8").Value = "ABCD0RGR014"
$ws.Range("V8").Value = "ZXC0987RGR014"

# Vehicle year fix on row 7
$ws.Range("O7").Value = 2015

# Leave the sheet with the same selection it was saved with
$ws.Range("V7:V8").Select()
